$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.196.65"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.789.54"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'226.67"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'32.06"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "'0.0689"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "2.049.46"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "'11.09"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "1.801.08"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "34.156.74"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "'0.622"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'4.18"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'68.12"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'245.10"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "0.0₃0778"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'10.82"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").Value = "'4.11"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'161.04"
$ws.Range("D26").Value = "'7.17"
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("D27").Value = "'16.34"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "1.449.14"
$ws.Range("E35").Value = "  +4.02%  "
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("E37").Value = "  +7.49%  "
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "'80.49"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "'13.48"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("D46").Value = "'6.05"
$ws.Range("E46").Value = "  +3.44%  "
$ws.Range("D47").Value = "'1.07"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "0.0₆0136"
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("D49").Value = "1.950.70"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "'106.02"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("E51").Value = "  +0.04%  "
